# Update ~ SnowMountain Field (Temp) (#86)
#
# This script reproduces the cell / selection / active-sheet changes made
# to Assets/Editor/JsonUtility/JsonUtility/Stage1/SnowMountain.xlsx.
# The NPC occupying the "ring box" branch sheets is renamed from
# 사진사 (photographer) to 아이작 (Isaac), a handful of new "Hold" /
# disappear steps are inserted on the ring-box-success sheet, a row is
# split off at the bottom of ring-box-0 / shifted up on ring-box-1 and
# ring-box-fail, and the active sheet/selection bookkeeping is updated to
# match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "1-4, branch"  (sheet13.xml)
# ---------------------------------------------------------------------
$ws13 = $wb.Worksheets.Item("1-4, branch")
$ws13.Range("A5").Value = "아이작"
$ws13.Range("A6").Value = "아이작"
$ws13.Range("A7").Value = "아이작"
$ws13.Range("A8").Value = "아이작"
$ws13.Range("A8").Select()

# ---------------------------------------------------------------------
# "1-4, string branch"  (sheet14.xml)
# ---------------------------------------------------------------------
$ws14 = $wb.Worksheets.Item("1-4, string branch")
$ws14.Range("A5").Value = "아이작"
$ws14.Range("A6").Value = "아이작"
$ws14.Range("A7").Value = "아이작"
$ws14.Range("A8").Value = "아이작"
$ws14.Range("A8").Select()

# ---------------------------------------------------------------------
# "1-4, ring box-0"  (sheet15.xml)
#   row 9's C9/F9 pair is split off into a new row 10; row 9 keeps
#   only a C9 = 1 marker.
# ---------------------------------------------------------------------
$ws15 = $wb.Worksheets.Item("1-4, ring box-0")
$ws15.Range("A5").Value = "아이작"
$ws15.Range("A6").Value = "아이작"
$ws15.Range("A7").Value = "아이작"
$ws15.Range("A8").Value = "아이작"
$ws15.Range("C10").Value = 11
$ws15.Range("F10").Value = "Reset"
$ws15.Range("F9").ClearContents()
$ws15.Range("C9").Value = 1
$ws15.Range("A8").Select()

# ---------------------------------------------------------------------
# "1-4, ring box-1"  (sheet16.xml)
#   row 9's C9/F9 pair moves up to row 7.
# ---------------------------------------------------------------------
$ws16 = $wb.Worksheets.Item("1-4, ring box-1")
$ws16.Range("A5").Value = "아이작"
$ws16.Range("C7").Value = 11
$ws16.Range("F7").Value = "Reset"
$ws16.Range("C9").ClearContents()
$ws16.Range("F9").ClearContents()
$ws16.Range("A5").Select()

# ---------------------------------------------------------------------
# "1-4, ring box-fail"  (sheet17.xml)
#   row 9's C9/F9 pair moves up to row 6.
# ---------------------------------------------------------------------
$ws17 = $wb.Worksheets.Item("1-4, ring box-fail")
$ws17.Range("A5").Value = "아이작"
$ws17.Range("C6").Value = 11
$ws17.Range("F6").Value = "Reset"
$ws17.Range("C9").ClearContents()
$ws17.Range("F9").ClearContents()
$ws17.Range("A5").Select()

# ---------------------------------------------------------------------
# "1-4, ring box-success"  (sheet18.xml)
#   Rename the NPC, add a handful of new Hold/remove rows describing
#   the ring/note popups, and append the ring-box "disappear" step plus
#   its trailing reset marker.
# ---------------------------------------------------------------------
$ws18 = $wb.Worksheets.Item("1-4, ring box-success")

$ws18.Range("A5").Value = "아이작"
$ws18.Range("A6").Value = "아이작"
$ws18.Range("A7").Value = "아이작"
$ws18.Range("C6").Value = 1
$ws18.Range("C7").Value = 1

$ws18.Range("C8").Value = 6
$ws18.Range("F8").Value = "Hold, 0.5, name=Stage 1/4 SnowMountain/Ring/Get Ring"
$ws18.Range("C9").Value = 6
$ws18.Range("F9").Value = "Hold, 0.5, name=Stage 1/4 SnowMountain/Ring/Remove Ring"

$ws18.Range("A10").Value = "아이작"
$ws18.Range("C10").Value = 1
$ws18.Range("A11").Value = "아이작"
$ws18.Range("C11").Value = 1
$ws18.Range("A12").Value = "아이작"
$ws18.Range("C12").Value = 1
$ws18.Range("A13").Value = "아이작"
$ws18.Range("C13").Value = 1

$ws18.Range("C14").Value = 6
$ws18.Range("F14").Value = "Hold, 0.5, name=Stage 1/4 SnowMountain/Ring/Get Note"
$ws18.Range("C15").Value = 6
$ws18.Range("F15").Value = "Hold, 0.5, name=Stage 1/4 SnowMountain/Ring/Remove Note"

$ws18.Range("A16").Value = "아이작"
$ws18.Range("C16").Value = 1
$ws18.Range("A17").Value = "아이작"
$ws18.Range("C17").Value = 1
$ws18.Range("A18").Value = "아이작"
$ws18.Range("C18").Value = 1
$ws18.Range("C19").Value = 1

$ws18.Range("B20").Value = "결혼반지함이 사라진다."
$ws18.Range("C20").Value = 6
$ws18.Range("F20").Value = "Hold, 0.5, name=Stage 1/4 SnowMountain/Ring/Disappear"

$ws18.Range("C21").Value = 11
$ws18.Range("F21").Value = "Reset"

# Make this the active sheet/selection, matching activeTab="17" and the
# tabSelected flag moving off of "1-4 main-2" (sheet4.xml) onto this sheet.
$ws18.Range("F21").Select()
$ws18.Activate()
